# Auto-generated edit script: update recalculated cost/profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled-runner refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 122
$ws.Range("I4").Value = 122
$ws.Range("K4").Value = 122
$ws.Range("M4").Value = -8
$ws.Range("H12").Value = 966.6667
$ws.Range("I12").Value = 966.6667
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 966.6667
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -796.6667
$ws.Range("N12").ClearContents()
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5350
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 3938.4688
$ws.Range("I132").Value = 1154.2222
$ws.Range("K132").Value = 3462.6666
$ws.Range("M132").Value = -932.6665999999996
$ws.Range("H137").Value = 40201
$ws.Range("I137").Value = 1820.3334
$ws.Range("K137").Value = 5461.0002
$ws.Range("M137").Value = -2911.0002
$ws.Range("H138").Value = 1889.7084
$ws.Range("J138").Value = 2590.2
$ws.Range("L138").Value = 7770.599999999999
$ws.Range("N138").Value = -18050.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4312.3125
$ws.Range("I2").Value = 4088.7273
$ws.Range("K2").Value = 4088.7273
$ws.Range("M2").Value = -3975.7273
$ws.Range("H32").Value = 51518.883
$ws.Range("I32").Value = 29869.555
$ws.Range("J32").Value = 181414.83
$ws.Range("K32").Value = 29869.555
$ws.Range("L32").Value = 181414.83
$ws.Range("M32").Value = -29582.555
$ws.Range("N32").Value = -181988.83
$ws.Range("H74").Value = 2985.9412
$ws.Range("I74").Value = 2647.5833
$ws.Range("K74").Value = 2647.5833
$ws.Range("M74").Value = -1773.5833
$ws.Range("H77").Value = 2985.9412
$ws.Range("I77").Value = 2647.5833
$ws.Range("K77").Value = 13237.9165
$ws.Range("M77").Value = -8869.916499999999
$ws.Range("H116").Value = 4312.3125
$ws.Range("I116").Value = 4088.7273
$ws.Range("K116").Value = 4088.7273
$ws.Range("M116").Value = -1794.7273

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4312.3125
$ws.Range("I3").Value = 4088.7273
$ws.Range("K3").Value = 4088.7273
$ws.Range("M3").Value = -3974.7273
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H94").Value = 629.4737
$ws.Range("I94").Value = 560.4706
$ws.Range("K94").Value = 560.4706
$ws.Range("M94").Value = -109.4706

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 1901.7142
$ws.Range("I58").Value = 1885.3334
$ws.Range("K58").Value = 1885.3334
$ws.Range("M58").Value = -1682.3334
$ws.Range("H60").Value = 9205.666999999999
$ws.Range("J60").Value = 9076.385
$ws.Range("L60").Value = 9076.385
$ws.Range("N60").Value = -10098.385
$ws.Range("H69").Value = 14500
$ws.Range("H72").Value = 14500
$ws.Range("H105").Value = 979.1667
$ws.Range("I105").Value = 785
$ws.Range("K105").Value = 785
$ws.Range("M105").Value = 962
$ws.Range("H132").Value = 1958.5927
$ws.Range("I132").Value = 2002.36
$ws.Range("J132").Value = 1411.5
$ws.Range("K132").Value = 6007.08
$ws.Range("L132").Value = 4234.5
$ws.Range("M132").Value = -3477.08
$ws.Range("N132").Value = -9294.5
$ws.Range("H134").Value = 3947.1667
$ws.Range("I134").Value = 4136.7
$ws.Range("K134").Value = 12410.1
$ws.Range("M134").Value = -9875.099999999999
$ws.Range("H136").Value = 1901.7142
$ws.Range("I136").Value = 1885.3334
$ws.Range("K136").Value = 5656.0002
$ws.Range("M136").Value = -3106.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1949.75
$ws.Range("I3").Value = 1949.75
$ws.Range("K3").Value = 5849.25
$ws.Range("M3").Value = -5737.25
$ws.Range("H5").Value = 546.2727
$ws.Range("I5").Value = 439.25
$ws.Range("K5").Value = 1317.75
$ws.Range("M5").Value = -1205.75
$ws.Range("H75").Value = 6406.5
$ws.Range("J75").Value = 7253.4546
$ws.Range("L75").Value = 21760.3638
$ws.Range("N75").Value = -23756.3638
$ws.Range("H78").Value = 6406.5
$ws.Range("J78").Value = 7253.4546
$ws.Range("L78").Value = 65281.0914
$ws.Range("N78").Value = -75265.0914
$ws.Range("H98").Value = 1666
$ws.Range("J98").Value = 1997.6666
$ws.Range("L98").Value = 5992.9998
$ws.Range("N98").Value = -8988.9998
$ws.Range("H107").Value = 882.36365
$ws.Range("I107").Value = 620.5714
$ws.Range("J107").Value = 952.8461
$ws.Range("K107").Value = 1861.7142
$ws.Range("L107").Value = 2858.5383
$ws.Range("M107").Value = 58.28579999999988
$ws.Range("N107").Value = -6698.5383
$ws.Range("H122").Value = 1022.89746
$ws.Range("I122").Value = 813.75
$ws.Range("J122").Value = 1076.871
$ws.Range("K122").Value = 7323.75
$ws.Range("L122").Value = 9691.839
$ws.Range("M122").Value = -4873.75
$ws.Range("N122").Value = -14591.839
$ws.Range("H132").Value = 976.1667
$ws.Range("I132").Value = 861
$ws.Range("J132").Value = 1157.1428
$ws.Range("K132").Value = 7749
$ws.Range("L132").Value = 10414.2852
$ws.Range("M132").Value = -5219
$ws.Range("N132").Value = -15474.2852
$ws.Range("H133").Value = 11412.333
$ws.Range("H135").Value = 546.2727
$ws.Range("I135").Value = 439.25
$ws.Range("K135").Value = 3953.25
$ws.Range("M135").Value = -1418.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4138
$ws.Range("I80").Value = 3700.6667
$ws.Range("K80").Value = 3700.6667
$ws.Range("M80").Value = -2702.6667
$ws.Range("H83").Value = 4138
$ws.Range("I83").Value = 3700.6667
$ws.Range("K83").Value = 18503.3335
$ws.Range("M83").Value = -13511.3335
$ws.Range("H126").Value = 3471
$ws.Range("I126").Value = 3333
$ws.Range("J126").Value = 3609
$ws.Range("K126").Value = 9999
$ws.Range("L126").Value = 10827
$ws.Range("M126").Value = -7529
$ws.Range("N126").Value = -15767
$ws.Range("H132").Value = 2397.3076
$ws.Range("I132").Value = 2453.5652
$ws.Range("K132").Value = 7360.6956
$ws.Range("M132").Value = -4830.6956

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 669.5
$ws.Range("I16").Value = 558.6667
$ws.Range("J16").Value = 1002
$ws.Range("K16").Value = 558.6667
$ws.Range("L16").Value = 1002
$ws.Range("M16").Value = -388.6667
$ws.Range("N16").Value = -1342
$ws.Range("H46").Value = 2385.9092
$ws.Range("I46").Value = 2374.5
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 2374.5
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -2186.5
$ws.Range("N46").Value = -2876
$ws.Range("H61").Value = 941.5
$ws.Range("I61").Value = 999
$ws.Range("J61").Value = 884
$ws.Range("K61").Value = 999
$ws.Range("L61").Value = 884
$ws.Range("M61").Value = -797
$ws.Range("N61").Value = -1288
$ws.Range("H113").Value = 941.5
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 884
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 884
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -5224
$ws.Range("H132").Value = 3179.8696
$ws.Range("I132").Value = 2071
$ws.Range("J132").Value = 4621.4
$ws.Range("K132").Value = 6213
$ws.Range("L132").Value = 13864.2
$ws.Range("M132").Value = -3683
$ws.Range("N132").Value = -18924.2
$ws.Range("H136").Value = 3379.0454
$ws.Range("I136").Value = 3012
$ws.Range("J136").Value = 3819.5
$ws.Range("K136").Value = 9036
$ws.Range("L136").Value = 11458.5
$ws.Range("M136").Value = -6486
$ws.Range("N136").Value = -16558.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2365.182
$ws.Range("I96").Value = 2201.9
$ws.Range("K96").Value = 2201.9
$ws.Range("M96").Value = -828.9000000000001
$ws.Range("H100").Value = 2315.6365
$ws.Range("I100").Value = 2763.6
$ws.Range("J100").Value = 1355.7142
$ws.Range("K100").Value = 5527.2
$ws.Range("L100").Value = 2711.4284
$ws.Range("M100").Value = -4986.2
$ws.Range("N100").Value = -3793.4284

Write-Output "Applied $([int]217) cell updates across 8 sheets."